$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(4,3,3,0),
  @(3,0,3,3),
  @(5,2,5,0),
  @(6,2,5,0),
  @(3,1,3,2),
  @(3,0,4,3),
  @(3,0,3,3),
  @(4,0,4,2),
  @(6,2,5,0),
  @(7,0,6,2),
  @(4,0,3,3),
  @(5,2,5,0),
  @(5,2,3,1),
  @(6,0,7,2),
  @(4,2,5,0),
  @(3,1,4,2),
  @(6,2,4,1),
  @(7,3,5,0),
  @(3,1,4,2),
  @(5,1,5,2),
  @(4,3,3,0),
  @(6,2,5,0),
  @(3,0,3,3),
  @(4,1,3,2),
  @(3,1,3,2),
  @(5,2,5,1),
  @(4,2,5,1),
  @(4,1,2,2),
  @(7,0,6,2),
  @(4,0,3,2),
  @(5,2,3,1),
  @(5,0,6,3),
  @(4,2,7,0),
  @(5,3,4,0),
  @(6,0,5,3),
  @(5,0,7,3),
  @(4,2,3,1),
  @(2,1,3,2),
  @(3,0,3,3)
)

$startRow = 2117
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$lastRow = $startRow + $data.Count
$ws.Range("A" + $lastRow).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2134
